$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.261.51"
$ws.Range("E2").Value = "  -2.58%  "

$ws.Range("D3").Value = "2.965.62"
$ws.Range("E3").Value = "  -2.59%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'587.52"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'141.55"
$ws.Range("E6").Value = "  -5.77%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  -2.55%  "

$ws.Range("D9").Value = "2.965.29"
$ws.Range("E9").Value = "  -2.60%  "

$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -6.46%  "

$ws.Range("D11").Value = "'5.77"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("E12").Value = "  +1.30%  "

$ws.Range("E13").Value = "  -3.37%  "

$ws.Range("D14").Value = "'33.83"
$ws.Range("E14").Value = "  -5.35%  "

$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("D16").Value = "3.459.15"
$ws.Range("E16").Value = "  -2.59%  "

$ws.Range("E17").Value = "  -1.28%  "

$ws.Range("D18").Value = "61.294.18"
$ws.Range("E18").Value = "  -2.54%  "

$ws.Range("D19").Value = "2.964.61"
$ws.Range("E19").Value = "  -2.80%  "

$ws.Range("D20").Value = "'445.48"
$ws.Range("E20").Value = "  -6.50%  "

$ws.Range("D21").Value = "'13.84"
$ws.Range("E21").Value = "  -2.46%  "

$ws.Range("E22").Value = "  -2.94%  "

$ws.Range("E23").Value = "  -2.17%  "

$ws.Range("D24").Value = "'81.15"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  -3.85%  "

$ws.Range("E26").Value = "  -8.97%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'9.85"
$ws.Range("E28").Value = "  -6.35%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("E30").Value = "  -0.17%  "

$ws.Range("D31").Value = "'6.83"
$ws.Range("E31").Value = "  -5.79%  "

$ws.Range("E32").Value = "  -6.31%  "

$ws.Range("D33").Value = "'26.93"
$ws.Range("E33").Value = "  -2.51%  "

$ws.Range("E34").Value = "  -3.75%  "

$ws.Range("E35").Value = "  -4.86%  "

$ws.Range("D36").Value = "0.0₃0777"
$ws.Range("E36").Value = "  -3.26%  "

$ws.Range("D37").Value = "'5.72"
$ws.Range("E37").Value = "  -2.29%  "

$ws.Range("D38").Value = "'50.11"
$ws.Range("E38").Value = "  -0.08%  "

$ws.Range("E39").Value = "  -5.67%  "

$ws.Range("D40").Value = "'9.09"
$ws.Range("E40").Value = "  +0.16%  "

$ws.Range("E41").Value = "  +4.81%  "

$ws.Range("D42").Value = "'2.76"
$ws.Range("E42").Value = "  -9.26%  "

$ws.Range("D43").Value = "'388.61"
$ws.Range("E43").Value = "  -8.14%  "

$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.681.70"
$ws.Range("E45").Value = "  -5.17%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.262"
$ws.Range("E46").Value = "  -7.62%  "

$ws.Range("D47").Value = "'37.14"
$ws.Range("E47").Value = "  -2.27%  "

$ws.Range("D48").Value = "'130.95"
$ws.Range("E48").Value = "  +2.61%  "

$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").Value = "'0.107"
$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("E51").Value = "  -1.24%  "
